$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 12186.5602457225
$ws.Range("C2").Value = 12125.3834500009
$ws.Range("E2").Value = 7459.6308145012
$ws.Range("F2").Value = -112.118572312411

# Row 3
$ws.Range("C3").Value = 12561.6711079526
$ws.Range("F3").Value = 378.491263663907

# Row 4
$ws.Range("C4").Value = 11675.5049615252
$ws.Range("F4").Value = 340.833401311903

# Row 5
$ws.Range("C5").Value = 7792.15554002028
$ws.Range("F5").Value = 153.335857315419

# Row 6
$ws.Range("C6").Value = 7663.24956666033
$ws.Range("F6").Value = 145.135470348352

# Row 7
$ws.Range("C7").Value = 11656.7677187487
$ws.Range("F7").Value = 328.967548054168

$wb.Save()
